$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.606.07"
$ws.Range("E2").Value = "  +6.16%  "
$ws.Range("D3").Value = "2.042.54"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.40%  "
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.17"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.83%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0755"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.21%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.912"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.57%  "
$ws.Range("D15").Value = "2.341.55"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +21.43%  "
$ws.Range("D18").Value = "2.026.45"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "37.431.78"
$ws.Range("E19").Value = "  +6.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.47%  "
$ws.Range("D21").Value = "0.0₃0873"
$ws.Range("E21").Value = "  +5.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("E24").Value = "  +26.20%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.122"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.114"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +29.31%  "
$ws.Range("E32").Value = "  +7.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0615"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.27%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +25.84%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  +21.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("E43").Value = "  +6.56%  "
$ws.Range("E44").Value = "  +6.89%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +23.52%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.39%  "
$ws.Range("D49").Value = "1.429.71"
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.16%  "
